$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 191-192, shifting existing rows 191:234 down to 193:236
$ws.Rows("191:192").Insert()

# Fill in the new row 191 (Primera) - week of 44543
$ws.Cells.Item(191, 1).Value = 1
$ws.Cells.Item(191, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(191, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(191, 4).Value = 44543
$ws.Cells.Item(191, 5).Value = 15
$ws.Cells.Item(191, 6).Value = 100112032
$ws.Cells.Item(191, 7).Value = "Zapallo italiano"
$ws.Cells.Item(191, 8).Value = "Huracán"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 120
$ws.Cells.Item(191, 11).Value = 5000
$ws.Cells.Item(191, 12).Value = 6000
$ws.Cells.Item(191, 13).Value = 5500
$ws.Cells.Item(191, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(191, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(191, 16).Value = 79
$ws.Cells.Item(191, 17).Value = 70
$ws.Cells.Item(191, 18).Value = "Hortaliza"

# Fill in the new row 192 (Segunda) - week of 44543
$ws.Cells.Item(192, 1).Value = 1
$ws.Cells.Item(192, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(192, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(192, 4).Value = 44543
$ws.Cells.Item(192, 5).Value = 15
$ws.Cells.Item(192, 6).Value = 100112032
$ws.Cells.Item(192, 7).Value = "Zapallo italiano"
$ws.Cells.Item(192, 8).Value = "Huracán"
$ws.Cells.Item(192, 9).Value = "Segunda"
$ws.Cells.Item(192, 10).Value = 120
$ws.Cells.Item(192, 11).Value = 4000
$ws.Cells.Item(192, 12).Value = 5000
$ws.Cells.Item(192, 13).Value = 4500
$ws.Cells.Item(192, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(192, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(192, 16).Value = 45
$ws.Cells.Item(192, 17).Value = 100
$ws.Cells.Item(192, 18).Value = "Hortaliza"
